$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","F","G","J","L","M","O")

$data = @(
    @(17.81455945005328, 10.46921337672148, 8.660633352591748, 36.37972241081521, 3.679568183624746, 10.76046124817535, 11.49398786328987, 17.00166471999923, 27.73487892589981),
    @(17.36159986530425, 10.26873412964327, 8.665931797615139, 36.51261448610473, 3.681819848806667, 10.78562023241519, 11.49438908246359, 16.89921508041572, 27.84915246853567),
    @(17.07930389074161, 10.14314948307899, 8.669954534369289, 36.60324849658839, 3.683275959311089, 10.80186041778952, 11.495792197518, 16.83788706638478, 27.92593984151888),
    @(16.96337661184304, 10.09139494202268, 8.671787822066563, 36.64244875748039, 3.683887897508859, 10.80867831550731, 11.49665585679657, 16.81330940474305, 27.95889252475712),
    @(16.94407812759291, 10.08276769151529, 8.672103968878167, 36.6490946090091, 3.683990632184861, 10.80982251329391, 11.49681692359931, 16.80925381001066, 27.96446450170106),
    @(17.0777438371181, 10.14245377746419, 8.66997847261397, 36.60376799811591, 3.683284136890722, 10.80195155617526, 11.49580266204337, 16.83755390437219, 27.92637753238851),
    @(17.65933980760163, 10.40063316897265, 8.662300857423956, 36.42366400629503, 3.68032932222629, 10.76897196224926, 11.49388673235529, 16.96602359850663, 27.77290335066007),
    @(18.7600856272816, 10.88515509783542, 8.65332837691708, 36.14247045419422, 3.675116049074276, 10.71055861702508, 11.4992628784925, 17.22966605954552, 27.52467585998161),
    @(19.53637413377304, 11.22532703726878, 8.650414948775413, 35.98013305643525, 3.671636327882476, 10.67141845895297, 11.50871827609574, 17.42941228660615, 27.3747003912483),
    @(19.88095401226867, 11.37615836726584, 8.6498816099044, 35.9159648321387, 3.670128607638212, 10.65442405101675, 11.51420077741347, 17.52135501759934, 27.31355922485147),
    @(20.01009280955014, 11.43267372973336, 8.649792920844, 35.89306355286539, 3.669568428914424, 10.648104642357, 11.51644534936424, 17.55630552027512, 27.2914294274782),
    @(19.98234206161719, 11.42052945577547, 8.649806993701665, 35.89793349396724, 3.669688595639235, 10.64946048968557, 11.51595447077763, 17.54877269179027, 27.29614990586295),
    @(19.89160594342006, 11.38082017909997, 8.6498720467987, 35.91405268292937, 3.670082306018687, 10.65390182777485, 11.51438207153443, 17.52422785669682, 27.31171806486475),
    @(19.83584888780355, 11.35641771531258, 8.649926626836873, 35.92410835976034, 3.670324864935985, 10.65663736615077, 11.51344082877159, 17.50921023008933, 27.32138737131684),
    @(19.51367271992818, 11.21538775385555, 8.650465688796768, 35.98452189485232, 3.671736369929658, 10.67254534733748, 11.50838363338154, 17.42342348773542, 27.37883902952657),
    @(19.31375084720966, 11.12783978323948, 8.650998805688578, 36.02406717813852, 3.672621509746807, 10.6825116009092, 11.50558277458233, 17.37105702715975, 27.4159011944555),
    @(19.19796184840418, 11.07711817431826, 8.651380024808224, 36.04772336312318, 3.673137702286416, 10.68832026171166, 11.50408306140325, 17.34103988901303, 27.43788483752441),
    @(19.15862410073121, 11.05988299893252, 8.65152192728589, 36.05588918437203, 3.673313694486081, 10.69030010242304, 11.50359443560453, 17.33089487491461, 27.44544245031748),
    @(19.33511645749137, 11.13719761342334, 8.650934339656727, 36.01976320881221, 3.672526552380747, 10.68144277988581, 11.5058694249277, 17.37662105321889, 27.41188685440444),
    @(19.9182947454342, 11.39250036479263, 8.649849869881928, 35.90928010435729, 3.669966372087084, 10.65259415526235, 11.51483936252448, 17.53143380784587, 27.30711751939659),
    @(20.29154263084394, 11.55583531308604, 8.649801057716529, 35.84522314689074, 3.668355851585448, 10.6344158171682, 11.52168300357834, 17.6333828114987, 27.24460972031119),
    @(20.09309100442464, 11.46899446303334, 8.649766932540071, 35.87866395700141, 3.669209697105037, 10.64405627843047, 11.51794111274294, 17.57890727305596, 27.27742406884251),
    @(19.32545970765373, 11.13296814750141, 8.650963251962821, 36.02170616466123, 3.672569459818805, 10.6819257479646, 11.50573948589166, 17.37410527808136, 27.41369963179226),
    @(18.46744572521197, 10.75668039812425, 8.655107412985956, 36.21079464135124, 3.676464559057673, 10.72569498501179, 11.49683688706679, 17.22966605954552, 27.58615635312848)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $addr = "{0}{1}" -f $cols[$j], $rowNum
        $ws.Range($addr).Value = $data[$i][$j]
    }
}
